$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E).
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the existing text-formatted Price/Volume cells)
# instead of auto-converting numeric- or percent-looking strings into
# actual numbers.
$updates = @{
    'D2' = "'270.08"
    'E2' = "'3.26%"
    'E3' = "'-1.89%"
    'D4' = "'4.709"
    'E4' = "'0.01%"
    'D5' = "'0.06095"
    'E5' = "'-1.73%"
    'D6' = "'6.743"
    'E6' = "'0.37%"
    'D7' = "'0.8583"
    'E7' = "'0.95%"
    'D8' = "'0.8924"
    'E8' = "'-2.54%"
    'D9' = "'0.1420"
    'E9' = "'0.69%"
    'D10' = "'0.05022"
    'E10' = "'7.93%"
    'D11' = "'0.07108"
    'E11' = "'0.24%"
    'D12' = "'0.03181"
    'E12' = "'1.10%"
    'D13' = "'0.09027"
    'E13' = "'-0.21%"
    'D14' = "'0.001530"
    'E14' = "'0.07%"
    'D15' = "'0.0006054"
    'E15' = "'-1.87%"
    'D16' = "'0.006081"
    'E16' = "'0.22%"
    'D17' = "'3.463"
    'E17' = "'-0.12%"
    'D18' = "'3.168"
    'E18' = "'-0.03%"
    'E19' = "'2.88%"
    'D21' = "'0.1300"
    'E21' = "'-0.71%"
    'D22' = "'3.833"
    'E22' = "'-6.15%"
    'E23' = "'-0.22%"
    'D24' = "'0.001185"
    'E24' = "'-2.03%"
    'D25' = "'0.004149"
    'E25' = "'9.12%"
    'E26' = "'-0.09%"
    'E27' = "'4.98%"
    'D40' = "'0.03954"
    'E40' = "'1.01%"
    'D41' = "'0.1119"
    'E41' = "'0.65%"
    'D42' = "'0.004191"
    'E42' = "'1.45%"
    'E43' = "'-7.99%"
    'D44' = "'0.01269"
    'E44' = "'-8.77%"
    'D45' = "'0.00005124"
    'E45' = "'-0.96%"
    'E47' = "'-31.84%"
    'E48' = "'24.13%"
    'D49' = "'0.00002100"
    'D50' = "'0.0002000"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
